$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-09 Monday", 2) | Out-Null
$d.Content.Find.Execute("822÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "267÷9=", 2) | Out-Null
$d.Content.Find.Execute("815÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "976÷8=", 2) | Out-Null
$d.Content.Find.Execute("712÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "208÷4=", 2) | Out-Null
$d.Content.Find.Execute("740÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "941÷2=", 2) | Out-Null
$d.Content.Find.Execute("590÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "189÷6=", 2) | Out-Null
$d.Content.Find.Execute("635÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "551÷3=", 2) | Out-Null
$d.Content.Find.Execute("485÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "954÷4=", 2) | Out-Null
$d.Content.Find.Execute("904÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "191÷5=", 2) | Out-Null
$d.Content.Find.Execute("622÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "789÷6=", 2) | Out-Null
$d.Content.Find.Execute("688÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "793÷4=", 2) | Out-Null
$d.Content.Find.Execute("261÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "616÷7=", 2) | Out-Null
$d.Content.Find.Execute("414÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "567÷9=", 2) | Out-Null
$d.Content.Find.Execute("131÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "530÷8=", 2) | Out-Null
$d.Content.Find.Execute("433÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "489÷2=", 2) | Out-Null
$d.Content.Find.Execute("377÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "705÷3=", 2) | Out-Null
$d.Content.Find.Execute("524÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "439÷8=", 2) | Out-Null
$d.Content.Find.Execute("431÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "245÷5=", 2) | Out-Null
$d.Content.Find.Execute("604÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "714÷5=", 2) | Out-Null
$d.Content.Find.Execute("573÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "458÷7=", 2) | Out-Null
$d.Content.Find.Execute("225÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "650÷6=", 2) | Out-Null
$d.Content.Find.Execute("323÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "932÷9=", 2) | Out-Null
$d.Content.Find.Execute("265÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "681÷7=", 2) | Out-Null
$d.Content.Find.Execute("262÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "238÷4=", 2) | Out-Null
$d.Content.Find.Execute("453÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "788÷2=", 2) | Out-Null
$d.Content.Find.Execute("773÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "968÷3=", 2) | Out-Null
